# Update the multiplication-practice table.
#
# Row 1 originally read [91×83=, 82×44=, 36×66=, 83×11=, 61×23=].
# The canonical edit drops the first cell (91×83=) and appends a new
# cell (71×79=) at the end, shifting everything else one slot left.
# Since every <w:tc> in the row shares identical formatting, writing
# the five resulting values into the five existing cells (in order)
# produces the exact same final content as delete-then-append.
#
# All the other populated rows (5, 10, 15, 20) keep their five cells
# and simply get new values in place.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-RowValues($table, [int]$rowIndex, [string[]]$values) {
    $row = $table.Rows.Item($rowIndex)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row.Cells.Item($i + 1).Range.Text = $values[$i]
    }
}

Set-RowValues $t 1  @("82×44=", "54×90=", "69×68=", "63×55=", "71×79=")
Set-RowValues $t 5  @("90×57=", "56×66=", "37×18=", "78×29=", "99×53=")
Set-RowValues $t 10 @("91×69=", "71×39=", "86×49=", "70×58=", "75×76=")
Set-RowValues $t 15 @("31×87=", "51×47=", "45×67=", "19×64=", "49×60=")
Set-RowValues $t 20 @("76×14=", "35×39=", "60×97=", "28×75=", "25×39=")

Write-Host "done"
